$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header columns: "<Name>_old" -> "<Name>_FV2404", "<Name>_new" -> "<Name>_FV2410"
$oldHeaders = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)
$newHeaders = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $oldHeaders[$i]
}
$ws.Cells.Item(1, 11).Value = "diff"
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, 12 + $i).Value = $newHeaders[$i]
}

# 2) Freeze the header row (split below row 1)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 3) Turn the data range into a real table ("Table1") spanning A1:U72
$rng = $ws.Range("A1:U72")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
